$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assay")

# Remove the example rows 3-5, keeping only a single example row (row 2).
# Delete from the bottom up so row indices stay valid.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Update the remaining example row (row 2) with the new example values.
$ws.Range("E2").Value2 = "libraryconstruction.txt"
$ws.Range("I2").Value2 = "single-end"
$ws.Range("J2").Value2 = "DPBO"
$ws.Range("K2").Value2 = "http://purl.obolibrary.org/obo/DPBO_0000086"
$ws.Range("L2").Value2 = "RNA-Seq"
$ws.Range("M2").Value2 = "EFO"
$ws.Range("N2").Value2 = "http://purl.obolibrary.org/obo/EFO_0008896"
$ws.Range("O2").Value2 = "Illumina HiSeq 2000"
$ws.Range("P2").Value2 = "OBI"
$ws.Range("Q2").Value2 = "http://purl.obolibrary.org/obo/OBI_0002001"
